{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet titlePara = null;\nlet datePara = null;\nlet firstPara = null;\n\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  if (t === \"Debate 1\") titlePara = p;\n  else if (t === \"F2025\") datePara = p;\n  else if (t === \"Debate stuff\") firstPara = p;\n}\n\nif (titlePara) titlePara.insertText(\"Debate 5\", \"Replace\");\nif (datePara) datePara.insertText(\"S2026\", \"Replace\");\n\nlet overviewPara = null;\nif (firstPara) {\n  // Add the new \"Overview\" heading immediately before the paragraph that\n  // used to read \"Debate stuff\", then update that paragraph's own text.\n  overviewPara = firstPara.insertParagraph(\"Overview\", \"Before\");\n  overviewPara.style = \"Heading 2\";\n  firstPara.insertText(\"Nothing to see here yet\", \"Replace\");\n}\nawait context.sync();\n\nif (overviewPara && firstPara) {\n  // Wrap a bookmark named \"overview\" around the new heading and the\n  // paragraph beneath it.\n  const bookmarkRange = overviewPara.getRange(\"Whole\").expandTo(firstPara.getRange(\"Whole\"));\n  bookmarkRange.insertBookmark(\"overview\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Title: \"Debate 1\" -> \"Debate 5\"\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"Debate 1`r\") {\n        $p.Range.Text = \"Debate 5\"\n        break\n    }\n}\n\n# Date: \"F2025\" -> \"S2026\"\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"F2025`r\") {\n        $p.Range.Text = \"S2026\"\n        break\n    }\n}\n\n# Find the \"Debate stuff\" paragraph, insert a new \"Overview\" heading right\n# before it, then update its own text.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"Debate stuff`r\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.InsertBefore(\"Overview`r\")\n\n    $overviewPara = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -eq \"Overview`r\") {\n            $overviewPara = $p\n            break\n        }\n    }\n    $overviewPara.Style = $d.Styles(\"Heading 2\")\n\n    $finalPara = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -eq \"Debate stuff`r\") {\n            $finalPara = $p\n            break\n        }\n    }\n    $finalPara.Range.Text = \"Nothing to see here yet\"\n\n    $bmRange = $d.Range($overviewPara.Range.Start, $finalPara.Range.End)\n    $d.Bookmarks.Add(\"overview\", $bmRange)\n}\n"}
